$d = $word.ActiveDocument

# Locate the "{{ postal_code }}" merge-field placeholder text.
$seeker = $d.Content
$found = $seeker.Find.Execute("{{ postal_code }}")

if (-not $found) {
    throw "Could not find '{{ postal_code }}' in the document"
}

# Re-acquire a plain Range for the matched span; using the Find range
# object directly for InsertXML appends after it instead of replacing.
$target = $d.Range($seeker.Start, $seeker.End)

# Re-split the run layout: merge "{{" with the following space into one
# run, break "postal_code" into "postal" / "_" / "code" runs, and merge
# the trailing space with "}}" into the final run. Character formatting
# (sz/szCs/lang) is preserved on every run.
$rPr = '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr>'

$runs = ''
$runs += '<w:r>' + $rPr + '<w:t xml:space="preserve">{{ </w:t></w:r>'
$runs += '<w:r>' + $rPr + '<w:t>postal</w:t></w:r>'
$runs += '<w:r>' + $rPr + '<w:t>_</w:t></w:r>'
$runs += '<w:r>' + $rPr + '<w:t>code</w:t></w:r>'
$runs += '<w:r>' + $rPr + '<w:t xml:space="preserve"> }}</w:t></w:r>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $runs + '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
